# Add 2022-Q3 data:
#  - new worksheet "2022-Q3" inserted right after "总计" (before "2022-Q2")
#    holding the per-fund breakdown for the new quarter
#  - the "总计" (summary) sheet gets a new row 2 for 2022-Q3 and all the
#    other quarter rows shift down by one

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q3" worksheet by duplicating "2022-Q2" (index 2)
#    so that it inherits the exact same column widths / header style /
#    per-row styling, then overwrite its name + values.
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item(2)
$srcSheet.Copy($srcSheet)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The source sheet only has 9 data rows (rows 2-10); 2022-Q3 needs 24
# (rows 2-25). Extend the column-A styling (bold/centered/bordered,
# matching the rest of column A) down to row 25 before writing values.
$q3.Range("A10").Copy()
$q3.Range("A11:A25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$q3.Range("A1").Select()

# row data: index, fund code, fund name, fund size, stock position,
#           position ratio, held value (100M), position rank
$q3Data = @(
  @(0,'008099','广发价值领先混合A','58.59','93.49','5.31','3.1111',4),
  @(1,'010761','华商甄选回报混合A','24.94','85.36','4.26','1.0624',2),
  @(2,'012671','嘉实核心蓝筹混合A','9.20','93.58','4.67','0.4296',7),
  @(3,'012420','广发价值领先混合C','7.55','93.49','5.31','0.4009',4),
  @(4,'012528','广发鑫睿一年持有期混合A','7.26','92.69','5.32','0.3862',6),
  @(5,'007368','浙商沪港深精选混合A','6.59','84.00','5.41','0.3565',5),
  @(6,'013607','广发睿恒进取一年持有期混合A','6.95','91.91','5.11','0.3551',7),
  @(7,'007177','浙商智能行业优选混合A','10.73','93.72','3.22','0.3455',6),
  @(8,'014734','广发睿合混合A','5.96','86.96','5.44','0.3242',5),
  @(9,'012529','广发鑫睿一年持有期混合C','4.74','92.69','5.32','0.2522',6),
  @(10,'014872','嘉实品质蓝筹一年持有期混合A','2.47','93.17','4.57','0.1129',7),
  @(11,'016049','华商甄选回报混合C','2.59','85.36','4.26','0.1103',2),
  @(12,'009126','嘉实基础产业优选股票A','1.74','93.85','4.81','0.0837',7),
  @(13,'014735','广发睿合混合C','1.47','86.96','5.44','0.0800',5),
  @(14,'005335','浙商全景消费混合A','1.61','91.41','4.33','0.0697',10),
  @(15,'008488','华商恒益稳健混合','2.39','52.53','2.35','0.0562',5),
  @(16,'014373','浙商全景消费混合C','0.98','91.41','4.33','0.0424',10),
  @(17,'007217','浙商智能行业优选混合C','0.90','93.72','3.22','0.0290',6),
  @(18,'013608','广发睿恒进取一年持有期混合C','0.47','91.91','5.11','0.0240',7),
  @(19,'012672','嘉实核心蓝筹混合C','0.42','93.58','4.67','0.0196',7),
  @(20,'007369','浙商沪港深精选混合C','0.32','84.00','5.41','0.0173',5),
  @(21,'009127','嘉实基础产业优选股票C','0.21','93.85','4.81','0.0101',7),
  @(22,'001900','诺安精选价值混合','0.12','85.83','3.00','0.0036',7),
  @(23,'014873','嘉实品质蓝筹一年持有期混合C','0.07','93.17','4.57','0.0032',7)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    # Text columns: prefix with "'" so numeric-looking strings (fund
    # codes, percentages, ...) are stored as text, not coerced to numbers.
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = "'" + $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert the 2022-Q3 row at the top
#    of the data (row 2) and push the existing rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Extend column-A styling down to the new row 9 (copy from row 8).
$total.Range("A8:D8").Copy()
$total.Range("A9:D9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$total.Range("A1").Select()

$totalData = @(
  @(0,'2022-Q3',24,7.69),
  @(1,'2022-Q2',9,5.37),
  @(2,'2022-Q1',3,4.19),
  @(3,'2021-Q4',3,1.74),
  @(4,'2021-Q3',3,4.36),
  @(5,'2021-Q2',3,1.76),
  @(6,'2021-Q1',7,2.36),
  @(7,'2020-Q4',1,0.02)
)

$r = 2
foreach ($row in $totalData) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = "'" + $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$total.Select()
